$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2455593528023
$ws.Range("K2").Value = 0.170354349596659
$ws.Range("L2").Value = 0.260087165113423
$ws.Range("N2").Value = 0.217833100525281

$ws.Range("B3").Value = 0.376242148576814
$ws.Range("K3").Value = 0.332091953956457
$ws.Range("L3").Value = 0.368533078351345
$ws.Range("N3").Value = 0.338459395559222

$ws.Range("B4").Value = 0.536095833727937
$ws.Range("K4").Value = 0.452298918345791
$ws.Range("L4").Value = 0.525032772881734
$ws.Range("N4").Value = 0.574624521327971

$ws.Range("B5").Value = 0.16455410482696
$ws.Range("K5").Value = 0.222495975817485
$ws.Range("L5").Value = 0.117168875472553
$ws.Range("N5").Value = 0.199940105738939
